$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cells whose values are removed in the updated dataset
$ws.Range("D4").ClearContents()
$ws.Range("AN4").ClearContents()
$ws.Range("AP4").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("E6").ClearContents()
$ws.Range("T6").ClearContents()
$ws.Range("AN6").ClearContents()
$ws.Range("AP6").ClearContents()
$ws.Range("AN7").ClearContents()
$ws.Range("AP7").ClearContents()
$ws.Range("AN8").ClearContents()
$ws.Range("AP8").ClearContents()

# Write refreshed capital-structure data for rows 2-9
# Row 2
$ws.Range("A2").Value = "Philippines"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "7"
$ws.Range("B2").NumberFormat = "General"
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Value = "Financial Svcs. (Non-bank & Insurance)"
$ws.Range("D2").Value = 0.0304
$ws.Range("E2").Value = -0.00506
$ws.Range("G2").Value = 0.05357421889532937
$ws.Range("H2").Value = 0.05357421889532937
$ws.Range("I2").Value = 0.04740574285863103
$ws.Range("J2").Value = 0.04019162993415164
$ws.Range("K2").Value = 57.138
$ws.Range("L2").Value = 0.2125780360584257
$ws.Range("M2").Value = 38.194
$ws.Range("N2").Value = 0.03103689257272875
$ws.Range("O2").Value = 0.6684518184045645
$ws.Range("P2").Value = 25.194
$ws.Range("Q2").Value = 0.02047294002925402
$ws.Range("R2").Value = 0.4409324792607372
$ws.Range("S2").Value = 13
$ws.Range("T2").Value = 0.3403675970047651
$ws.Range("U2").Value = 121.332
$ws.Range("V2").Value = 0.09859580692345198
$ws.Range("W2").Value = 0.02401574803149606
$ws.Range("X2").Value = 0.0200779880695942
$ws.Range("Y2").Value = 0.003937759961901866
$ws.Range("Z2").Value = 0.309650588113314
$ws.Range("AA2").Value = -0
$ws.Range("AB2").Value = 0.0200779880695942
$ws.Range("AC2").Value = -0.0200779880695942
$ws.Range("AD2").Value = 268.93
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 268.93
$ws.Range("AG2").Value = 147.598
$ws.Range("AH2").Value = 0.1793428607630391
$ws.Range("AI2").Value = 0.240293287728271
$ws.Range("AJ2").Value = 0.1070949166955691
$ws.Range("AK2").Value = 0.1479172053291002
$ws.Range("AL2").Value = 0.167
$ws.Range("AM2").Value = -3.203
$ws.Range("AN2").Value = 17.02088607594937
$ws.Range("AO2").Value = 76.29940119760478
$ws.Range("AP2").Value = 9.341645569620253
$ws.Range("AQ2").Value = -3.978145488604433

# Row 3
$ws.Range("A3").Value = "Philippines"
$ws.Range("B3").Value = "The Philippine Stock Exchange, Inc. (PSE:PSE)"
$ws.Range("C3").Value = "Financial Svcs. (Non-bank & Insurance)"
$ws.Range("D3").Value = 0.00006
$ws.Range("E3").Value = -0.0108
$ws.Range("G3").Value = 0.549618320610687
$ws.Range("H3").Value = 0.549618320610687
$ws.Range("I3").Value = 0.4847328244274809
$ws.Range("J3").Value = 0.3557756013250756
$ws.Range("K3").Value = 11.7
$ws.Range("L3").Value = 0.4465648854961832
$ws.Range("M3").Value = 13
$ws.Range("N3").Value = 0.04965622612681436
$ws.Range("O3").Value = 1.111111111111111
$ws.Range("P3").Value = -0
$ws.Range("Q3").Value = -0
$ws.Range("R3").Value = -0
$ws.Range("S3").Value = 13
$ws.Range("T3").Value = 1
$ws.Range("U3").Value = 42.2
$ws.Range("V3").Value = 0.1611917494270435
$ws.Range("W3").Value = 0.1134820562560621
$ws.Range("X3").Value = 0.02011631620898607
$ws.Range("Y3").Value = 0.093365740047076
$ws.Range("Z3").Value = 0.7660818713450295
$ws.Range("AA3").Value = 0.2725532384420171
$ws.Range("AB3").Value = 0.02014208982312409
$ws.Range("AC3").Value = 0.252411148618893
$ws.Range("AD3").Value = 1.33
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 1.33
$ws.Range("AG3").Value = -40.87
$ws.Range("AH3").Value = 0.005054535780792764
$ws.Range("AI3").Value = 0.01226597805035507
$ws.Range("AJ3").Value = -0.1849907210428643
$ws.Range("AK3").Value = -0.6170919522874831
$ws.Range("AL3").Value = 0.167
$ws.Range("AM3").Value = -3.203
$ws.Range("AN3").Value = 0.08417721518987342
$ws.Range("AO3").Value = 76.04790419161675
$ws.Range("AP3").Value = -2.586708860759494
$ws.Range("AQ3").Value = -3.965032781767093

# Row 4
$ws.Range("A4").Value = "Philippines"
$ws.Range("B4").Value = "Ferronoux Holdings, Inc. (PSE:FERRO)"
$ws.Range("C4").Value = "Financial Svcs. (Non-bank & Insurance)"
$ws.Range("G4").Value = -0
$ws.Range("H4").Value = -0
$ws.Range("I4").Value = -0
$ws.Range("J4").Value = -0
$ws.Range("K4").Value = -0.084
$ws.Range("L4").Value = 3.5
$ws.Range("M4").Value = -0
$ws.Range("N4").Value = -0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = -0
$ws.Range("Q4").Value = -0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("U4").Value = 0.012
$ws.Range("V4").Value = 0.0005405405405405405
$ws.Range("W4").Value = -0.03243243243243243
$ws.Range("X4").Value = 0.0200779880695942
$ws.Range("Y4").Value = -0.05251042050202663
$ws.Range("Z4").Value = -0.00927357032457496
$ws.Range("AA4").Value = 0
$ws.Range("AB4").Value = 0.0200779880695942
$ws.Range("AC4").Value = -0.0200779880695942
$ws.Range("AD4").Value = 0
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 0
$ws.Range("AG4").Value = -0.012
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = -0.0005408328826392645
$ws.Range("AK4").Value = -0.004497751124437781
$ws.Range("AL4").Value = 0
$ws.Range("AM4").Value = 0

# Row 5
$ws.Range("A5").Value = "Philippines"
$ws.Range("B5").Value = "Citystate Savings Bank, Inc. (PSE:CSB)"
$ws.Range("C5").Value = "Financial Svcs. (Non-bank & Insurance)"
$ws.Range("D5").Value = 0.0505
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0.305
$ws.Range("L5").Value = 0.04951298701298701
$ws.Range("M5").Value = -0
$ws.Range("N5").Value = -0
$ws.Range("O5").Value = -0
$ws.Range("P5").Value = -0
$ws.Range("Q5").Value = -0
$ws.Range("R5").Value = -0
$ws.Range("S5").Value = 0
$ws.Range("U5").Value = 4.13
$ws.Range("V5").Value = 0.2333333333333333
$ws.Range("W5").Value = 0.02401574803149606
$ws.Range("X5").Value = 0.0200779880695942
$ws.Range("Y5").Value = 0.003937759961901866
$ws.Range("Z5").Value = 1.271938880858972
$ws.Range("AA5").Value = 0
$ws.Range("AB5").Value = 0.0200779880695942
$ws.Range("AC5").Value = -0.0200779880695942
$ws.Range("AD5").Value = 0
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 0
$ws.Range("AG5").Value = -4.13
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = -0.3043478260869565
$ws.Range("AK5").Value = -0.4503816793893129
$ws.Range("AL5").Value = 0
$ws.Range("AM5").Value = 0

# Row 6
$ws.Range("A6").Value = "Philippines"
$ws.Range("B6").Value = "Prime Media Holdings, Inc. (PSE:PRIM)"
$ws.Range("C6").Value = "Financial Svcs. (Non-bank & Insurance)"
$ws.Range("D6").Value = 0.0133
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = -0.015
$ws.Range("L6").Value = -0.2631578947368421
$ws.Range("M6").Value = -0
$ws.Range("N6").Value = -0
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = -0
$ws.Range("Q6").Value = -0
$ws.Range("R6").Value = 0
$ws.Range("S6").Value = 0
$ws.Range("U6").Value = 0.196
$ws.Range("V6").Value = 0.01568
$ws.Range("W6").Value = 0.004934210526315789
$ws.Range("X6").Value = 0.0200779880695942
$ws.Range("Y6").Value = -0.01514377754327841
$ws.Range("Z6").Value = -0.02030637691485572
$ws.Range("AA6").Value = -0
$ws.Range("AB6").Value = 0.0200779880695942
$ws.Range("AC6").Value = -0.0200779880695942
$ws.Range("AD6").Value = 0
$ws.Range("AE6").Value = 0
$ws.Range("AF6").Value = 0
$ws.Range("AG6").Value = -0.196
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = -0
$ws.Range("AJ6").Value = -0.0159297789336801
$ws.Range("AK6").Value = 0.06210392902408111
$ws.Range("AL6").Value = 0
$ws.Range("AM6").Value = 0

# Row 7
$ws.Range("A7").Value = "Philippines"
$ws.Range("B7").Value = "Philippine Savings Bank (PSE:PSB)"
$ws.Range("C7").Value = "Financial Svcs. (Non-bank & Insurance)"
$ws.Range("D7").Value = 0.0475
$ws.Range("E7").Value = 0.00068
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 45
$ws.Range("L7").Value = 0.1923899102180419
$ws.Range("M7").Value = 25.1
$ws.Range("N7").Value = 0.05174190888476603
$ws.Range("O7").Value = 0.5577777777777778
$ws.Range("P7").Value = 25.1
$ws.Range("Q7").Value = 0.05174190888476603
$ws.Range("R7").Value = 0.5577777777777778
$ws.Range("S7").Value = 0
$ws.Range("T7").Value = 0
$ws.Range("U7").Value = 72.2
$ws.Range("V7").Value = 0.1488352916924345
$ws.Range("W7").Value = 0.06767935027823734
$ws.Range("X7").Value = 0.02400503278470608
$ws.Range("Y7").Value = 0.04367431749353125
$ws.Range("Z7").Value = 0.2891581159599456
$ws.Range("AA7").Value = 0
$ws.Range("AB7").Value = 0.02441938114446802
$ws.Range("AC7").Value = -0.02441938114446802
$ws.Range("AD7").Value = 252.5
$ws.Range("AE7").Value = 0
$ws.Range("AF7").Value = 252.5
$ws.Range("AG7").Value = 180.3
$ws.Range("AH7").Value = 0.3423264642082429
$ws.Range("AI7").Value = 0.259800390986727
$ws.Range("AJ7").Value = 0.2709648331830478
$ws.Range("AK7").Value = 0.2004001333777926
$ws.Range("AL7").Value = 0
$ws.Range("AM7").Value = 0

# Row 8
$ws.Range("A8").Value = "Philippines"
$ws.Range("B8").Value = "Makati Finance Corporation (PSE:MFIN)"
$ws.Range("C8").Value = "Financial Svcs. (Non-bank & Insurance)"
$ws.Range("D8").Value = -0.0488
$ws.Range("E8").Value = -0.235
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0.206
$ws.Range("L8").Value = 0.08442622950819671
$ws.Range("M8").Value = 0.067
$ws.Range("N8").Value = 0.004466666666666667
$ws.Range("O8").Value = 0.325242718446602
$ws.Range("P8").Value = 0.067
$ws.Range("Q8").Value = 0.004466666666666667
$ws.Range("R8").Value = 0.325242718446602
$ws.Range("S8").Value = 0
$ws.Range("T8").Value = 0
$ws.Range("U8").Value = 1.6
$ws.Range("V8").Value = 0.1066666666666667
$ws.Range("W8").Value = 0.0206
$ws.Range("X8").Value = 0.02767287699596824
$ws.Range("Y8").Value = -0.007072876995968237
$ws.Range("Z8").Value = 0.117195004803074
$ws.Range("AA8").Value = 0
$ws.Range("AB8").Value = 0.02875774191106082
$ws.Range("AC8").Value = -0.02875774191106082
$ws.Range("AD8").Value = 15.1
$ws.Range("AE8").Value = 0
$ws.Range("AF8").Value = 15.1
$ws.Range("AG8").Value = 13.5
$ws.Range("AH8").Value = 0.5016611295681063
$ws.Range("AI8").Value = 0.5852713178294574
$ws.Range("AJ8").Value = 0.4736842105263158
$ws.Range("AK8").Value = 0.5578512396694215
$ws.Range("AL8").Value = 0
$ws.Range("AM8").Value = 0

# Row 9
$ws.Range("A9").Value = "Philippines"
$ws.Range("B9").Value = "GMA Holdings, Inc. (PSE:GMAP)"
$ws.Range("C9").Value = "Financial Svcs. (Non-bank & Insurance)"
$ws.Range("D9").Value = 0.467
$ws.Range("E9").Value = 0.09050000000000001
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0.7924528301886793
$ws.Range("J9").Value = 0.6293007769145395
$ws.Range("K9").Value = 0.026
$ws.Range("L9").Value = 0.4905660377358491
$ws.Range("M9").Value = 0.027
$ws.Range("N9").Value = 0.00006485707422531828
$ws.Range("O9").Value = 1.038461538461539
$ws.Range("P9").Value = 0.027
$ws.Range("Q9").Value = 0.00006485707422531828
$ws.Range("R9").Value = 1.038461538461539
$ws.Range("S9").Value = 0
$ws.Range("T9").Value = 0
$ws.Range("U9").Value = 0.994
$ws.Range("V9").Value = 0.002387701177035791
$ws.Range("W9").Value = 1.368421052631579
$ws.Range("X9").Value = 0.0200779880695942
$ws.Range("Y9").Value = 1.348343064561985
$ws.Range("Z9").Value = -0.1031128404669261
$ws.Range("AA9").Value = -0.06488899061570154
$ws.Range("AB9").Value = 0.0200779880695942
$ws.Range("AC9").Value = -0.08496697868529574
$ws.Range("AD9").Value = 0
$ws.Range("AE9").Value = 0
$ws.Range("AF9").Value = 0
$ws.Range("AG9").Value = -0.994
$ws.Range("AH9").Value = 0
$ws.Range("AI9").Value = 0
$ws.Range("AJ9").Value = -0.002393415939090695
$ws.Range("AK9").Value = 1.024742268041237
$ws.Range("AL9").Value = 0
$ws.Range("AM9").Value = 0

